# amazon.xlsx update - "all works at this moment"
#
# The test workflow on the "Amazon" sheet has most of its rows (7-19)
# flipped from "not runnable" (RUN = n) to "runnable" (RUN = y) -- except
# row 10, which is intentionally left as a negative test and stays "n".
# The "Sign out" step's element locator is also switched from an id-based
# locator to an xpath locator, and the sheet's last selection / one row's
# height are nudged as a side effect of the edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Amazon")

# --- Flip RUN column to "y" for the now-working steps (rows 7-19), ---
# --- but keep row 10 ("this is negative test") as "n".             ---
$runRows = 7..19
foreach ($r in $runRows) {
    if ($r -eq 10) {
        $ws.Cells.Item($r, 1).Value = "n"
    } else {
        $ws.Cells.Item($r, 1).Value = "y"
    }
}

# --- Update the "Sign out" locator from an id-based locator to an xpath one ---
$newSignOutLocator = "<@!{xpath=//span[contains(text(),'Sign Out')]}>"
$ws.Range("C12").Value = $newSignOutLocator
$ws.Range("C19").Value = $newSignOutLocator

# --- Row 19 height nudges slightly after the content edits ---
$ws.Rows.Item(19).RowHeight = 13.8

# --- Update the sheet's remembered selection / active cell ---
$ws.Range("C23").Select()
